$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values per row (B,C,D,E,G) - F column remains 0 / unchanged
$ws.Range("B2").Value = 0.3048080303191223
$ws.Range("C2").Value = 0.3127903958511391
$ws.Range("D2").Value = 3.900430680208489
$ws.Range("E2").Value = 8.660232485948974
$ws.Range("G2").Value = 13.17826159232772

$ws.Range("B3").Value = 0.01514828764759746
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 26.21740644021617
$ws.Range("E3").Value = 8.660232485948974
$ws.Range("G3").Value = 36.56058179708087

$ws.Range("B4").Value = 0.003994804209775715
$ws.Range("C4").Value = 0.00007097389502863649
$ws.Range("D4").Value = 0.1575252929769615
$ws.Range("E4").Value = 8.660232485948974
$ws.Range("G4").Value = 8.82182355703074

$ws.Range("B5").Value = 0.04763786555579896
$ws.Range("C5").Value = 0.002777888934908601
$ws.Range("D5").Value = 0.8054896365839992
$ws.Range("E5").Value = 8.660232485948974
$ws.Range("G5").Value = 9.516137877023681
